$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a "Price"-column (D) text value. Many of these values look
# like plain numbers (e.g. "22.03"); prefixing with a leading apostrophe
# forces Excel to store them as literal text (matching the source
# <c t="inlineStr"><is><t>...</t></is></c> cells), exactly like typing
# '22.03 into the cell in the Excel UI.
function Set-PriceCell($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
}

# Row 2 - Bitcoin
Set-PriceCell "D2" "27.141.01"
$ws.Range("E2").Value = "  +0.79%  "

# Row 3 - Ethereum
Set-PriceCell "D3" "1.572.19"
$ws.Range("E3").Value = "  +1.28%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +1.11%  "

# Row 5 - BNB
Set-PriceCell "D5" "210.87"
$ws.Range("E5").Value = "  +2.17%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.68%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.59%  "

# Row 8 - Solana
Set-PriceCell "D8" "22.03"
$ws.Range("E8").Value = "  +0.37%  "

# Row 9 - Cardano
Set-PriceCell "D9" "0.248"
$ws.Range("E9").Value = "  +0.13%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.75%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.72%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-PriceCell "D12" "1.791.47"
$ws.Range("E12").Value = "  +1.02%  "

# Row 13 - WrappedEther
Set-PriceCell "D13" "1.576.78"
$ws.Range("E13").Value = "  +1.64%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.83%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.43%  "

# Row 16 - WrappedBTC
Set-PriceCell "D16" "27.156.36"
$ws.Range("E16").Value = "  +0.88%  "

# Row 17 - Litecoin
Set-PriceCell "D17" "62.17"
$ws.Range("E17").Value = "  +0.85%  "

# Row 18 - ShibaInu
Set-PriceCell "D18" "0.0₃0705"
$ws.Range("E18").Value = "  -0.90%  "

# Row 19 - BitcoinCash
Set-PriceCell "D19" "216.47"
$ws.Range("E19").Value = "  -0.23%  "

# Row 20 - Chainlink
Set-PriceCell "D20" "7.41"
$ws.Range("E20").Value = "  +1.58%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.78%  "

# Row 22 - Uniswap
Set-PriceCell "D22" "4.15"
$ws.Range("E22").Value = "  +1.70%  "

# Row 23 - Avalanche
Set-PriceCell "D23" "9.22"
$ws.Range("E23").Value = "  +0.32%  "

# Row 24 - Toncoin
Set-PriceCell "D24" "1.95"
$ws.Range("E24").Value = "  +0.42%  "

# Row 25 - Monero
Set-PriceCell "D25" "154.25"
$ws.Range("E25").Value = "  +0.34%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  -0.18%  "

# Row 27 - EthereumClassic
Set-PriceCell "D27" "15.16"
$ws.Range("E27").Value = "  +1.15%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  +1.62%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.55%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +5.54%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +1.20%  "

# Row 32 - Filecoin
Set-PriceCell "D32" "3.25"
$ws.Range("E32").Value = "  +1.02%  "

# Row 33 - InternetComputer(DFINITY)
Set-PriceCell "D33" "3.19"
$ws.Range("E33").Value = "  +2.59%  "

# Row 34 - Maker
Set-PriceCell "D34" "1.436.44"
$ws.Range("E34").Value = "  +1.81%  "

# Row 35 - TrustWalletToken
Set-PriceCell "D35" "1.11"
$ws.Range("E35").Value = "  +15.04%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  +1.15%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  +2.80%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +0.86%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  +1.07%  "

# Rows 40-42 were re-ranked: MXToken / FraxShare / ARBITRUM rotate order.
# Row 40 becomes FraxShare
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-PriceCell "D40" "5.88"
$ws.Range("E40").Value = "  +4.08%  "

# Row 41 becomes ARBITRUM
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-PriceCell "D41" "0.812"
$ws.Range("E41").Value = "  +0.64%  "

# Row 42 becomes MXToken
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-PriceCell "D42" "2.38"
$ws.Range("E42").Value = "  +3.42%  "

# Row 43 - PaxDollar
$ws.Range("E43").Value = "  +0.78%  "

# Row 44 - WEMIXToken (unchanged)

# Row 45 - Aave
Set-PriceCell "D45" "64.74"
$ws.Range("E45").Value = "  +0.43%  "

# Row 46 - RenderToken
Set-PriceCell "D46" "1.74"
$ws.Range("E46").Value = "  +0.80%  "

# Row 47 - RocketPoolETH
Set-PriceCell "D47" "1.709.85"
$ws.Range("E47").Value = "  +1.34%  "

# Row 48 - Quant
Set-PriceCell "D48" "85.75"
$ws.Range("E48").Value = "  -1.80%  "

# Row 49 - BabyDogeCoin
$ws.Range("E49").Value = "  +2.01%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  -0.49%  "

# Row 51 - Algorand
Set-PriceCell "D51" "0.0962"
$ws.Range("E51").Value = "  +0.42%  "
